$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

foreach ($ws in @($ws1, $ws4)) {
    $ws.Range("F4").Value = 12297
    $ws.Range("F5").Value = 1268
    $ws.Range("F6").Value = 135
    $ws.Range("F15").Value = 39
    $ws.Range("F16").Value = 363
    $ws.Range("F17").Value = 3042
    $ws.Range("F18").Value = 90
    $ws.Range("F19").Value = 935
    $ws.Range("F20").Value = 13
    $ws.Range("F22").Value = 27
}
